# Reto 2 - Entrega Final: update cached data values and selection/active-sheet state.

$wb = $excel.ActiveWorkbook

# --- Update data values on "Tablas datos" (columns B and F, rows 15-21) ---
$wsDatos = $wb.Worksheets.Item("Tablas datos")

$wsDatos.Range("B15").Value = 107.605
$wsDatos.Range("B16").Value = 232.732
$wsDatos.Range("B17").Value = 502.91
$wsDatos.Range("B18").Value = 798.71100000000001
$wsDatos.Range("B19").Value = 1400.3030000000001
$wsDatos.Range("B20").Value = 2435.2829999999999
$wsDatos.Range("B21").Value = 3297.1689999999999

$wsDatos.Range("F16").Value = 239.696
$wsDatos.Range("F17").Value = 493.37299999999999
$wsDatos.Range("F18").Value = 795.67
$wsDatos.Range("F19").Value = 1407.2070000000001
$wsDatos.Range("F20").Value = 2400.0329999999999
$wsDatos.Range("F21").Value = 3112.1680000000001

# --- Update selection / active sheet state for every sheet ---
$wsReq0 = $wb.Worksheets.Item("Requerimiento 0")
$wsReq1 = $wb.Worksheets.Item("Requerimiento 1")
$wsReq2 = $wb.Worksheets.Item("Requerimiento 2")
$wsReq3 = $wb.Worksheets.Item("Requerimiento 3")
$wsReq4 = $wb.Worksheets.Item("Requerimiento 4")

$wsReq0.Range("V20").Select()
$wsReq1.Range("U16").Select()
$wsReq2.Range("T13").Select()
$wsReq3.Range("V13").Select()
$wsReq4.Range("V16").Select()

# "Tablas datos" ends up as the active/selected tab with L20 selected.
$wsDatos.Activate()
$wsDatos.Range("L20").Select()
